$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    foreach ($cell in $t.Range.Cells) {
        $cr = $cell.Range
        if ($cr.Font.NameAscii -eq "Helvetica") {
            $cr.Font.NameOther = "Helvetica"
            $cr.Font.NameFarEast = "Helvetica"
            $cr.Font.NameBi = "Helvetica"
        }
    }
}
